# Populate the worksheet with the two text values that become shared
# strings "aaa" / "vvv" (A1, A2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "aaa"
$ws.Range("A2").Value = "vvv"

# Move / leave the active selection on A3, matching the saved selection
# in the worksheet's sheetView.
[void]$ws.Range("A3").Select()

# Page setup: paper size 9 (A4) and portrait orientation, as in the target.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
